$d = $word.ActiveDocument

# --- Hunk 1: paragraph with "Script_calculateKinematics ... WORKING HERE..." ---
$target1 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like 'Script_calculate*WORKING HERE*') {
        $target1 = $p
        break
    }
}
if ($null -eq $target1) { throw "Paragraph for hunk 1 not found" }

$xml11 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:cs="Courier"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>cript_calculate</w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Kinematics</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>, which will move the origin to the initial pellet location and calculate a variety of kinematic features (velocity, aperture, etc.)</w:t></w:r></w:p>
'@
$target1.Range.InsertXML($xml11)

# --- Hunk 2: the empty paragraph immediately preceding "script_analyze_interp_trajectories" ---
$target2b = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like 'script_analyze_interp_trajectories*') {
        $target2b = $p
        break
    }
}
if ($null -eq $target2b) { throw "Paragraph for script_analyze_interp_trajectories not found" }
$target2a = $target2b.Previous()
if ($null -eq $target2a) { throw "Preceding empty paragraph not found" }

$xml13 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/><w:color w:val="008000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/><w:color w:val="00B050"/></w:rPr><w:t>script_analyze_interp_trajectories</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>. This will identify individual reaches from each trial and calculate some summary statistics for each session (average trajectories, variances, etc.). see script comments for details</w:t></w:r></w:p>
'@
$xml14 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/><w:color w:val="008000"/></w:rPr></w:pPr></w:p>
'@

$target2a.Range.InsertXML($xml13)
$target2b.Range.InsertXML($xml14)

Write-Output "Done"
